# Update the "variables" sheet: column N (row 2-10) values change 1 -> 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 14).Value = 2
}

# Move the active selection/view: drop the Q1 frozen top-left scroll position
# and move the selected cell from S18 to M19
$null = $ws.Activate()
$null = $ws.Range("M19").Select()
